$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should carry the same bold/
# bordered/centered header style as the existing header cells (e.g. H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-39: column I is (mostly) a constant 1, column J mirrors
# column H -- except row 38, which breaks the pattern (I38=3, J38=5).
for ($r = 2; $r -le 39; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2

    if ($r -eq 38) {
        $ws.Cells.Item($r, 9).Value = 3
        $ws.Cells.Item($r, 10).Value = 5
    } else {
        $ws.Cells.Item($r, 9).Value = 1
        $ws.Cells.Item($r, 10).Value = $hVal
    }
}
